$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seasons")

# Australia's row (row 2) gets its generic season names replaced with the
# more specific Noongar six-season names.
$ws.Range("B2").Value = "Summer-Birak"
$ws.Range("C2").Value = "Summer-Bunuru"
$ws.Range("D2").Value = "Autumn-Bunuru"
$ws.Range("E2").Value = "Autumn-Djeran"
$ws.Range("F2").Value = "Autumn-Djeran"
$ws.Range("G2").Value = "Winter-Makuru"
$ws.Range("H2").Value = "Winter-Makuru"
$ws.Range("I2").Value = "Winter-Dijiba"
$ws.Range("J2").Value = "Spring-Dijiba"
$ws.Range("K2").Value = "Spring-Kambarang"
$ws.Range("L2").Value = "Spring-Kambarang"
$ws.Range("M2").Value = "Summer-Birak"

$ws.Range("L18").Select()
